$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.11624
$ws.Range("H2").Value = 6.34872
$ws.Range("I2").Value = 0.1897594766532197
$ws.Range("J2").Value = 0.1897594766532197
$ws.Range("M2").Value = 46.29121633333333
$ws.Range("N2").Value = 138.873649
$ws.Range("O2").Value = 0.3133663986859022
$ws.Range("P2").Value = 0.3133663986859022
$ws.Range("Q2").Value = 97.96332365325333
$ws.Range("R2").Value = 881.66991287928
$ws.Range("S2").Value = 0.059464243815341
$ws.Range("T2").Value = 0.05946424381534101
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.11624
$ws.Range("H3").Value = 6.34872
$ws.Range("I3").Value = 0.1897594766532197
$ws.Range("J3").Value = 0.1897594766532197
$ws.Range("M3").Value = 46.81622333333333
$ws.Range("O3").Value = 0.3169204109998198
$ws.Range("P3").Value = 0.3169204109998198
$ws.Range("Q3").Value = 99.07436446693333
$ws.Range("R3").Value = 891.6692802023999
$ws.Range("S3").Value = 0.0601386513320491
$ws.Range("T3").Value = 0.0601386513320491
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.11624
$ws.Range("H4").Value = 6.34872
$ws.Range("I4").Value = 0.1897594766532197
$ws.Range("J4").Value = 0.1897594766532197
$ws.Range("M4").Value = 38.53544233333333
$ws.Range("N4").Value = 115.606327
$ws.Range("O4").Value = 0.2608640200510233
$ws.Range("P4").Value = 0.2608640200510233
$ws.Range("Q4").Value = 81.55024448349332
$ws.Range("R4").Value = 733.95220035144
$ws.Range("S4").Value = 0.0495014199225372
$ws.Range("T4").Value = 0.0495014199225372
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.11624
$ws.Range("H5").Value = 6.34872
$ws.Range("I5").Value = 0.1897594766532197
$ws.Range("J5").Value = 0.1897594766532197
$ws.Range("M5").Value = 16.07945366666667
$ws.Range("N5").Value = 48.238361
$ws.Range("O5").Value = 0.1088491702632547
$ws.Range("P5").Value = 0.1088491702632547
$ws.Range("Q5").Value = 34.02798302754666
$ws.Range("R5").Value = 306.25184724792
$ws.Range("S5").Value = 0.02065516158329242
$ws.Range("T5").Value = 0.02065516158329243
$ws.Range("I6").Value = 0.6160274054778138
$ws.Range("J6").Value = 0.6160274054778138
$ws.Range("M6").Value = 46.29121633333333
$ws.Range("N6").Value = 138.873649
$ws.Range("O6").Value = 0.3133663986859022
$ws.Range("P6").Value = 0.3133663986859022
$ws.Range("Q6").Value = 318.024128051225
$ws.Range("R6").Value = 2862.217152461025
$ws.Range("S6").Value = 0.1930422895464025
$ws.Range("T6").Value = 0.1930422895464025
$ws.Range("I7").Value = 0.6160274054778138
$ws.Range("J7").Value = 0.6160274054778138
$ws.Range("M7").Value = 46.81622333333333
$ws.Range("O7").Value = 0.3169204109998198
$ws.Range("P7").Value = 0.3169204109998198
$ws.Range("Q7").Value = 321.63096551675
$ws.Range("R7").Value = 2894.67868965075
$ws.Range("S7").Value = 0.1952316585311814
$ws.Range("T7").Value = 0.1952316585311814
$ws.Range("I8").Value = 0.6160274054778138
$ws.Range("J8").Value = 0.6160274054778138
$ws.Range("M8").Value = 38.53544233333333
$ws.Range("N8").Value = 115.606327
$ws.Range("O8").Value = 0.2608640200510233
$ws.Range("P8").Value = 0.2608640200510233
$ws.Range("Q8").Value = 264.741378988175
$ws.Range("R8").Value = 2382.672410893575
$ws.Range("S8").Value = 0.1606993854545443
$ws.Range("T8").Value = 0.1606993854545443
$ws.Range("I9").Value = 0.6160274054778138
$ws.Range("J9").Value = 0.6160274054778138
$ws.Range("M9").Value = 16.07945366666667
$ws.Range("N9").Value = 48.238361
$ws.Range("O9").Value = 0.1088491702632547
$ws.Range("P9").Value = 0.1088491702632547
$ws.Range("Q9").Value = 110.467052649025
$ws.Range("R9").Value = 994.203473841225
$ws.Range("S9").Value = 0.06705407194568561
$ws.Range("T9").Value = 0.06705407194568561
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.793503666666667
$ws.Range("H10").Value = 5.380511
$ws.Range("I10").Value = 0.1608202836929164
$ws.Range("J10").Value = 0.1608202836929164
$ws.Range("M10").Value = 46.29121633333333
$ws.Range("N10").Value = 138.873649
$ws.Range("O10").Value = 0.3133663986859022
$ws.Range("P10").Value = 0.3133663986859022
$ws.Range("Q10").Value = 83.02346622829323
$ws.Range("R10").Value = 747.211196054639
$ws.Range("S10").Value = 0.05039567313649433
$ws.Range("T10").Value = 0.05039567313649433
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.793503666666667
$ws.Range("H11").Value = 5.380511
$ws.Range("I11").Value = 0.1608202836929164
$ws.Range("J11").Value = 0.1608202836929164
$ws.Range("M11").Value = 46.81622333333333
$ws.Range("O11").Value = 0.3169204109998198
$ws.Range("P11").Value = 0.3169204109998198
$ws.Range("Q11").Value = 83.9650682078189
$ws.Range("R11").Value = 755.68561387037
$ws.Range("S11").Value = 0.05096723040506667
$ws.Range("T11").Value = 0.05096723040506667
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.793503666666667
$ws.Range("H12").Value = 5.380511
$ws.Range("I12").Value = 0.1608202836929164
$ws.Range("J12").Value = 0.1608202836929164
$ws.Range("M12").Value = 38.53544233333333
$ws.Range("N12").Value = 115.606327
$ws.Range("O12").Value = 0.2608640200510233
$ws.Range("P12").Value = 0.2608640200510233
$ws.Range("Q12").Value = 69.11345712145523
$ws.Range("R12").Value = 622.021114093097
$ws.Range("S12").Value = 0.0419522257098802
$ws.Range("T12").Value = 0.0419522257098802
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.793503666666667
$ws.Range("H13").Value = 5.380511
$ws.Range("I13").Value = 0.1608202836929164
$ws.Range("J13").Value = 0.1608202836929164
$ws.Range("M13").Value = 16.07945366666667
$ws.Range("N13").Value = 48.238361
$ws.Range("O13").Value = 0.1088491702632547
$ws.Range("P13").Value = 0.1088491702632547
$ws.Range("Q13").Value = 28.83855910916344
$ws.Range("R13").Value = 259.547031982471
$ws.Range("S13").Value = 0.01750515444147518
$ws.Range("T13").Value = 0.01750515444147518
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.3724043333333333
$ws.Range("H14").Value = 1.117213
$ws.Range("I14").Value = 0.03339283417605023
$ws.Range("J14").Value = 0.03339283417605023
$ws.Range("M14").Value = 46.29121633333333
$ws.Range("N14").Value = 138.873649
$ws.Range("O14").Value = 0.3133663986859022
$ws.Range("P14").Value = 0.3133663986859022
$ws.Range("Q14").Value = 17.23904955780411
$ws.Range("R14").Value = 155.151446020237
$ws.Range("S14").Value = 0.01046419218766438
$ws.Range("T14").Value = 0.01046419218766438
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.3724043333333333
$ws.Range("H15").Value = 1.117213
$ws.Range("I15").Value = 0.03339283417605023
$ws.Range("J15").Value = 0.03339283417605023
$ws.Range("M15").Value = 46.81622333333333
$ws.Range("O15").Value = 0.3169204109998198
$ws.Range("P15").Value = 0.3169204109998198
$ws.Range("Q15").Value = 17.43456443963445
$ws.Range("R15").Value = 156.91107995671
$ws.Range("S15").Value = 0.01058287073152266
$ws.Range("T15").Value = 0.01058287073152266
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.3724043333333333
$ws.Range("H16").Value = 1.117213
$ws.Range("I16").Value = 0.03339283417605023
$ws.Range("J16").Value = 0.03339283417605023
$ws.Range("M16").Value = 38.53544233333333
$ws.Range("N16").Value = 115.606327
$ws.Range("O16").Value = 0.2608640200510233
$ws.Range("P16").Value = 0.2608640200510233
$ws.Range("Q16").Value = 14.35076571185011
$ws.Range("R16").Value = 129.156891406651
$ws.Range("S16").Value = 0.008710988964061663
$ws.Range("T16").Value = 0.008710988964061663
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.3724043333333333
$ws.Range("H17").Value = 1.117213
$ws.Range("I17").Value = 0.03339283417605023
$ws.Range("J17").Value = 0.03339283417605023
$ws.Range("M17").Value = 16.07945366666667
$ws.Range("N17").Value = 48.238361
$ws.Range("O17").Value = 0.1088491702632547
$ws.Range("P17").Value = 0.1088491702632547
$ws.Range("Q17").Value = 5.988058223099222
$ws.Range("R17").Value = 53.89252400789299
$ws.Range("S17").Value = 0.003634782292801523
$ws.Range("T17").Value = 0.003634782292801523
